# Apply the "Ajout draft mapping" change:
#  - Update the Date metadata value on the "Metadata" sheet
#  - Add a new "Mapping: Spécification métier vers l'extension RORContactDescription"
#    column (AL) on the "Elements" sheet, with a draft mapping value of
#    "description" for the Extension.value[x] row.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump the Date property ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-12T09:15:29+00:00"

# --- Elements sheet: add the new mapping column ---
$ws = $wb.Worksheets.Item("Elements")

# Header for the new column AL (column 38)
$ws.Cells.Item(1, 38).Value = "Mapping: Spécification métier vers l'extension RORContactDescription"

# Data rows 2-6: blank except for the Extension.value[x] row which gets the draft mapping
$ws.Cells.Item(2, 38).Value = ""
$ws.Cells.Item(3, 38).Value = ""
$ws.Cells.Item(4, 38).Value = ""
$ws.Cells.Item(5, 38).Value = ""
$ws.Cells.Item(6, 38).Value = "description"

# Match the bestFit column width used for the new column
$ws.Columns.Item(38).ColumnWidth = 74.24609375

# Copy the formatting (bold header row / wrapped body rows) from the existing
# "Mapping: RIM Mapping" column (AK) onto the new column so AL matches the
# same header/body styling.
$src = $ws.Range("AK1:AK6")
$dst = $ws.Range("AL1:AL6")
$src.Copy()
$dst.PasteSpecial(-4122)
